# Apply the "swap hotel_info / review_info sheet roles + add State column" edit.
#
# Original layout:
#   Worksheet #1 (physical sheet1.xml, r:id=rId1): name "hotel_info"
#       header (9 cols: STR,Hotel_Name,City,Zip,TA_ReviewURL,Tripadvisor_Hotel_Name,
#               English_Reviews_num,Local_Rank,Total_Reviews_num) + 1 data row
#   Worksheet #2 (physical sheet2.xml, r:id=rId2): name "review_info"
#       header only (25 cols: STR,reviewer_ID,...,response_text)
#
# Target layout:
#   Worksheet #1 (same physical sheet, rId1): renamed "review_info"
#       header only (25 cols, same review_info header list), no data row
#   Worksheet #2 (same physical sheet, rId2): renamed "hotel_info"
#       header (10 cols: STR,Hotel_Name,State,City,Zip,TA_ReviewURL,
#               Tripadvisor_Hotel_Name,English_Reviews_num,Local_Rank,Total_Reviews_num)
#       + 1 data row, with a new "State" = "Louisiana" value inserted after Hotel_Name

$wb = $excel.ActiveWorkbook

$sheetA = $wb.Worksheets.Item("hotel_info")    # will become "review_info"
$sheetB = $wb.Worksheets.Item("review_info")   # will become "hotel_info"

# review_info headers (in column order)
$reviewHeaders = @(
    "STR", "reviewer_ID", "reviewer_name", "Review_ID", "Date_of_scraping",
    "ReviewURL", "Tripadvisor_gcode", "Tripadvisor_dcode", "Tripadvisor_rcode",
    "review_date", "review_title", "review_content", "review_rating",
    "trip_month", "trip_purpose", "value", "rooms", "Location", "Cleanliness",
    "Sleep Quality", "Service", "Picture(yes=1)", "respondent", "response_date",
    "response_text"
)

# hotel_info headers (in column order, with the new "State" column inserted)
$hotelHeaders = @(
    "STR", "Hotel_Name", "State", "City", "Zip", "TA_ReviewURL",
    "Tripadvisor_Hotel_Name", "English_Reviews_num", "Local_Rank", "Total_Reviews_num"
)

# hotel_info data row (with "Louisiana" inserted between the hotel name and city)
$hotelData = @(
    34709, "Residence Inn New Orleans Downtown", "Louisiana", "New Orleans", 70130,
    "https://www.tripadvisor.com/Hotel_Review-g60864-d93150-Reviews-Residence_Inn_New_Orleans_Downtown-New_Orleans_Louisiana.html",
    "Residence Inn by Marriott New Orleans Downtown", "1111", "47", "1144"
)

# --- Rebuild sheetA (rId1) as the new "review_info": header-only, 25 columns ---
$sheetA.Cells.Clear()
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $sheetA.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- Rebuild sheetB (rId2) as the new "hotel_info": 10-column header + 1 data row ---
$sheetB.Cells.Clear()
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $sheetB.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}
# Columns H, I, J ("1111", "47", "1144") are numeric-looking but must stay text,
# matching the source data (stored as shared strings, not numbers).
$textCols = @(8, 9, 10)
for ($i = 0; $i -lt $hotelData.Length; $i++) {
    $col = $i + 1
    $cell = $sheetB.Cells.Item(2, $col)
    if ($textCols -contains $col) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $hotelData[$i]
}

# --- Rename the sheets to reflect their swapped roles ---
# Use a temporary name to avoid a name collision while swapping.
$sheetA.Name = "review_info_tmp_swap"
$sheetB.Name = "hotel_info"
$sheetA.Name = "review_info"
